$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the "Conversión del día" note with today's rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$newNote = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.81 = 41299.34 pesos`n✅ 41299.34 pesos = 9.77 = 972.2 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newNote

# --- Sheet "tasas": refresh the Binance/transfi rate inputs ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 101.89
$ws2.Range("O10").Value = 4207.99
$ws2.Range("N12").Value = 4227
$ws2.Range("O12").Value = 99.505
